$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A32").Value = "LORENZA SIMONCELLI"
$ws.Range("B32").Value = "Riccardo Versini | Modium"
$ws.Range("C32").Value = "Davide Simoncelli | Avanzi"
$ws.Range("D32").Value = "Carlo  Stedile | Mai una gioia"
$ws.Range("E32").Value = "Federico Fasanelli | SBARX"
$ws.Range("F32").Value = "Emanuele Miorandi | Rita Levi’s"
